$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'27.051.18"
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.32%  '
$c = $ws.Range("D3")
$c.Value = "'1.829.57"
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  -0.32%  '
$c = $ws.Range("D5")
$c.Value = "'313.04"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.17%  '
$c = $ws.Range("D6")
$c.Value = "'1.007"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '
$c = $ws.Range("D7")
$c.Value = "'0.4590"
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.83%  '
$c = $ws.Range("D8")
$c.Value = "'0.3706"
$c.Style = "Normal"
$ws.Range("E8").Value = '  +2.05%  '
$c = $ws.Range("D9")
$c.Value = "'0.07340"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.73%  '
$c = $ws.Range("D10")
$c.Value = "'0.8740"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.56%  '
$c = $ws.Range("D11")
$c.Value = "'0.07956"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +4.19%  '
$c = $ws.Range("D12")
$c.Value = "'19.81"
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.48%  '
$c = $ws.Range("D13")
$c.Value = "'1.788.97"
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.59%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D14")
$c.Value = "'6.586"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.85%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D15")
$c.Value = "'5.333"
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.14%  '
$c = $ws.Range("D16")
$c.Value = "'91.68"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.74%  '
$c = $ws.Range("D17")
$c.Value = "'1.009"
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.15%  '
$c = $ws.Range("D18")
$c.Value = "'0.000008886"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +3.17%  '
$c = $ws.Range("D19")
$c.Value = "'1.009"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("C20").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c = $ws.Range("D20")
$c.Value = "'27.386.32"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range("D21")
$c.Value = "'14.71"
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.59%  '
$c = $ws.Range("D22")
$c.Value = "'5.104"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.11%  '
$c = $ws.Range("D23")
$c.Value = "'10.55"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.12%  '
$c = $ws.Range("D24")
$c.Value = "'2.146.22"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +2.35%  '
$c = $ws.Range("D25")
$c.Value = "'152.81"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.16%  '
$c = $ws.Range("D26")
$c.Value = "'1.850"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.98%  '
$c = $ws.Range("D27")
$c.Value = "'18.43"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.84%  '
$c = $ws.Range("D28")
$c.Value = "'2.044"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.43%  '
$c = $ws.Range("D29")
$c.Value = "'5.147"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.50%  '
$ws.Range("E30").Value = '  -0.80%  '
$c = $ws.Range("D31")
$c.Value = "'0.08865"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.49%  '
$c = $ws.Range("D32")
$c.Value = "'2.967"
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.46%  '
$c = $ws.Range("D33")
$c.Value = "'0.7316"
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.58%  '
$c = $ws.Range("D34")
$c.Value = "'4.442"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.33%  '
$c = $ws.Range("D35")
$c.Value = "'1.134"
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.40%  '
$c = $ws.Range("D36")
$c.Value = "'1.073"
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.82%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D37")
$c.Value = "'0.05242"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D38")
$c.Value = "'0.01940"
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.52%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D39")
$c.Value = "'2.423"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -3.68%  '
$c = $ws.Range("D40")
$c.Value = "'2.948"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.74%  '
$c = $ws.Range("D41")
$c.Value = "'7.170"
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.08%  '
$c = $ws.Range("D42")
$c.Value = "'0.5151"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.09%  '
$c = $ws.Range("D43")
$c.Value = "'0.1630"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.05%  '
$c = $ws.Range("D44")
$c.Value = "'8.221"
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.82%  '
$c = $ws.Range("D45")
$c.Value = "'0.4830"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.42%  '
$c = $ws.Range("D46")
$c.Value = "'1.008"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.25%  '
$c = $ws.Range("D47")
$c.Value = "'10.17"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.01%  '
$c = $ws.Range("D48")
$c.Value = "'102.55"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.97%  '
$ws.Range("E49").Value = '  -0.49%  '
$c = $ws.Range("D51")
$c.Value = "'64.92"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.48%  '
